# AWS SA LAB Credentials.xlsx - "Add files via upload" edit
#
# The Password column (D2:D9) is changed from 8 distinct per-participant
# random passwords to a single shared password "intel@123" for every
# participant, matching the rest of the sheet's formatting (column E's
# black Verdana font over a bordered cell) and getting a hyperlink just
# like the Login Id column already has.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Replace each password in D2:D9 with the shared password "intel@123".
$passwordRange = $ws.Range("D2:D9")
$passwordRange.Value = "intel@123"

# 2) Match D2:D9's look to the rest of the data rows (same font/border as
#    column E : plain, non-underlined, black Verdana 10 on a thin border).
$passwordRange.Font.Name = "Verdana"
$passwordRange.Font.Size = 10
$passwordRange.Font.Underline = $false
$passwordRange.Font.Color = 0

# 3) Widen column C slightly to fit the (unchanged) Login Id content.
$ws.Columns.Item(3).ColumnWidth = 45.5703125

# 4) Hyperlink the new password cells, same as every other "account info"
#    column on the sheet (D2 gets its own link, D3:D9 share one).
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:intel@123", "", "", "intel@123") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3:D9"), "mailto:intel@123", "", "", "intel@123") | Out-Null

# 5) Leave the selection parked on the edited range, like the author did.
$ws.Range("D2:D9").Select()
